$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Laptop"
$ws.Range("A4").Value = "Pendriver"

$ws.Range("B11").Select()
